# Auto-generated script to apply cryptos.xlsx price/volume updates
# (commit: Updated cryptos list on Sun Sep 10 20:40:16 UTC 2023 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellRef, $text) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = '@'
    $c.Value = $text
    $c.Style = 'Normal'
}

$ws.Range('D2').Value = '25.958.78'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '1.627.70'
$ws.Range('E3').Value = '  -0.96%  '
Set-TextValue 'D4' '1.01'
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '214.21'
$ws.Range('E5').Value = '  -0.87%  '
$ws.Range('E6').Value = '  -0.78%  '
$ws.Range('E7').Value = '  +0.01%  '
Set-TextValue 'D8' '0.251'
$ws.Range('E8').Value = '  -1.84%  '
$ws.Range('E9').Value = '  -3.09%  '
Set-TextValue 'D10' '18.48'
$ws.Range('E10').Value = '  -5.62%  '
$ws.Range('E11').Value = '  -1.16%  '
$ws.Range('D12').Value = '1.856.86'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('D14').Value = '1.620.15'
$ws.Range('E14').Value = '  -1.86%  '
Set-TextValue 'D15' '0.528'
$ws.Range('E15').Value = '  -2.98%  '
$ws.Range('D16').Value = '25.980.16'
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('E17').Value = '  -3.00%  '
Set-TextValue 'D18' '61.38'
$ws.Range('E18').Value = '  -3.21%  '
Set-TextValue 'D20' '192.47'
$ws.Range('E20').Value = '  -0.98%  '
$ws.Range('E21').Value = '  -2.58%  '
$ws.Range('E22').Value = '  -3.40%  '
$ws.Range('E23').Value = '  -1.91%  '
$ws.Range('E24').Value = '  +0.32%  '
Set-TextValue 'D25' '143.88'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').Value = '  -3.73%  '
$ws.Range('E28').Value = '  -2.00%  '
Set-TextValue 'D29' '15.23'
$ws.Range('E29').Value = '  -1.82%  '
$ws.Range('E30').Value = '  -1.33%  '
$ws.Range('E31').Value = '  -2.07%  '
$ws.Range('E33').Value = '  -5.28%  '
$ws.Range('E34').Value = '  -2.55%  '
$ws.Range('E35').Value = '  -2.63%  '
$ws.Range('D36').Value = '1.128.35'
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('E37').Value = '  -5.65%  '
$ws.Range('E38').Value = '  -1.48%  '
Set-TextValue 'D39' '0.523'
$ws.Range('E39').Value = '  -3.24%  '
$ws.Range('E40').Value = '  -2.18%  '
$ws.Range('E41').Value = '  -1.00%  '
$ws.Range('D42').Value = '1.766.66'
$ws.Range('E42').Value = '  -0.78%  '
Set-TextValue 'D43' '0.763'
$ws.Range('E43').Value = '  -4.26%  '
$ws.Range('E44').Value = '  -5.12%  '
$ws.Range('E45').Value = '  +1.91%  '
Set-TextValue 'D46' '54.41'
$ws.Range('E46').Value = '  -3.60%  '
$ws.Range('E47').Value = '  -0.63%  '
$ws.Range('D48').Value = '0.0₇0980'
$ws.Range('E48').Value = '  -16.46%  '
$ws.Range('E49').Value = '  -0.46%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D50' '7.51'
$ws.Range('E50').Value = '  -2.81%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
Set-TextValue 'D51' '1.01'
$ws.Range('E51').Value = '  +0.05%  '
